$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-01 23:03:53"

$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-01 23:03:48"

$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-01 23:03:53"

$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null
$zhcn.Columns.Item(3).AutoFit() | Out-Null
$dede.Columns.Item(3).AutoFit() | Out-Null
